######################################################################
# "Generate Report for Archive"
#
# The localization-status report moves from "Ready for handoff" to
# "In Translation": update every cell that shows that status (the
# Overview sheet's per-language status columns, plus each language
# sheet's own "Status" column) and let the (now narrower) Status
# columns shrink to fit the new text, just like Excel's column
# AutoFit would after the content got shorter.
######################################################################

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# ---- Overview sheet: columns E (zh-cn) and F (de-de) ---------------
$overview = $wb.Worksheets.Item("Overview")
$overviewRows = @(2, 3, 4)
foreach ($r in $overviewRows) {
    if ($overview.Cells.Item($r, 5).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($overview.Cells.Item($r, 6).Value2 -eq $oldStatus) {
        $overview.Cells.Item($r, 6).Value = $newStatus
    }
}
# Status text got shorter -> columns resize to fit the new content.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# ---- Per-language sheets: column C ("Status") -----------------------
$langSheets = @("zh-cn", "de-de")
foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in @(2, 3, 4)) {
        if ($ws.Cells.Item($r, 3).Value2 -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
